# Make the near-white background of the "Shape 220" picture on slide 1
# transparent by adding a color-change (clrChange FCFCFC -> alpha 0) to its
# blip fill, matching PowerPoint's "Set Transparent Color" picture tool.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)

# 0xFCFCFC == RGB(252,252,252) == 16579836
$sh.PictureFormat.TransparencyColor = 16579836
